$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BVAbIC")

$ws.Range("B2").Value = 611545701.8008617
$ws.Range("C2").Value = 4703397286.103693
$ws.Range("D2").Value = 622395909.8230397
$ws.Range("E2").Value = 764011492.1200054
$ws.Range("F2").Value = 177148614.5522801
$ws.Range("G2").Value = 4115544.45443973
$ws.Range("H2").Value = 19185897.01236145
$ws.Range("I2").Value = 59928449.72486594
$ws.Range("J2").Value = 261806312.9550531
$ws.Range("K2").Value = 120946743.487877
$ws.Range("L2").Value = 40965767.46167467
$ws.Range("M2").Value = 120936501.2780334
$ws.Range("N2").Value = 8969708.49343178
$ws.Range("O2").Value = 63623067.8871683
$ws.Range("P2").Value = 662348983.1354555
$ws.Range("Q2").Value = 30802013.47168675
$ws.Range("R2").Value = 48541470.04478592
$ws.Range("S2").Value = 3518253.538589437
$ws.Range("T2").Value = 38184778.75455172
$ws.Range("U2").Value = 65187133.96902961
$ws.Range("V2").Value = 1088750217.882793
$ws.Range("W2").Value = 1955239287.88846
$ws.Range("X2").Value = 5459611826.547406
$ws.Range("Y2").Value = 1317057209.917042
$ws.Range("Z2").Value = 1730199002.765774
$ws.Range("AA2").Value = 761645129.3893812
$ws.Range("AB2").Value = 655023870.2898594
$ws.Range("AC2").Value = 90009033.85006481
$ws.Range("AD2").Value = 1896269756.484216
$ws.Range("AE2").Value = 9287419501.511457
$ws.Range("AF2").Value = 7562938457.122041
$ws.Range("AG2").Value = 14336626334.75562
$ws.Range("AH2").Value = 425444857.3307816
$ws.Range("AI2").Value = 4416688274.298625
$ws.Range("AJ2").Value = 406456720.3863209
